$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
# (Hoja1 is also the workbook's ActiveSheet)

# Update the date in A1 (45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update price values in column D for rows 33-36
$ws.Range("D33").Value = 1560
$ws.Range("D34").Value = 2015
$ws.Range("D35").Value = 2249
$ws.Range("D36").Value = 2516
